# Add a new "frostedheart:wild_rubber_dandelion" plant-temperature row,
# placed right after the existing "frostedheart:rubber_dandelion" row
# (row 95) with identical min/max fertilize/grow/survive numbers and
# survive_snow / survive_blizzard / dead / will_die flags.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 96

$ws.Cells.Item($newRow, 1).Value = "frostedheart:wild_rubber_dandelion"
$ws.Cells.Item($newRow, 2).Value = 14
$ws.Cells.Item($newRow, 3).Value = 22
$ws.Cells.Item($newRow, 4).Value = 6
$ws.Cells.Item($newRow, 5).Value = 26
$ws.Cells.Item($newRow, 6).Value = -20
$ws.Cells.Item($newRow, 7).Value = 30
$ws.Cells.Item($newRow, 8).Value = $true
$ws.Cells.Item($newRow, 9).Value = $true
$ws.Cells.Item($newRow, 10).Value = "minecraft:dead_bush"
$ws.Cells.Item($newRow, 11).Value = $true

# Columns A-G inherit the column's own default cell style automatically.
# H/I (survive_snow / survive_blizzard) use the same "AR ADGothicJP Medium"
# font as the rest of the data rows, but are explicitly styled on each row
# (not through a column default) - match that font explicitly.
$ws.Cells.Item($newRow, 8).Font.Name = "AR ADGothicJP Medium"
$ws.Cells.Item($newRow, 9).Font.Name = "AR ADGothicJP Medium"

# J/K (dead / will_die) use the "DengXian" font family used elsewhere for
# these two trailing columns.
$ws.Cells.Item($newRow, 10).Font.Name = "等线"
$ws.Cells.Item($newRow, 10).Font.Family = 4
$ws.Cells.Item($newRow, 11).Font.Name = "等线"
$ws.Cells.Item($newRow, 11).Font.Family = 4

# Leave the view pointed at the freshly-entered row, same as someone typing
# the new plant's data in directly below the last existing row.
$ws.Rows($newRow).Select() | Out-Null
